$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.954332188745184
$ws.Cells.Item(2, 4).Value = 0.2236207094243312
$ws.Cells.Item(2, 5).Value = 0.2705100878955164
$ws.Cells.Item(2, 6).Value = 1.281511570911228
$ws.Cells.Item(2, 7).Value = 0.002446657893996875
$ws.Cells.Item(2, 10).Value = 0.3910156462670926
$ws.Cells.Item(2, 12).Value = 0.912800330224087
$ws.Cells.Item(2, 15).Value = 3.122830585786858
$ws.Cells.Item(3, 2).Value = 1.823371147136641
$ws.Cells.Item(3, 4).Value = 0.2233165457572142
$ws.Cells.Item(3, 5).Value = 0.2641648254123226
$ws.Cells.Item(3, 6).Value = 1.301250591704658
$ws.Cells.Item(3, 7).Value = 0.002450243460326706
$ws.Cells.Item(3, 10).Value = 0.3756362661057722
$ws.Cells.Item(3, 12).Value = 0.8099460157775127
$ws.Cells.Item(3, 15).Value = 3.150675514364764
$ws.Cells.Item(4, 2).Value = 1.742978581418981
$ws.Cells.Item(4, 4).Value = 0.2231998219204669
$ws.Cells.Item(4, 5).Value = 0.2603610044638174
$ws.Cells.Item(4, 6).Value = 1.314526286040419
$ws.Cells.Item(4, 7).Value = 0.00245256200940658
$ws.Cells.Item(4, 10).Value = 0.3663186081013947
$ws.Cells.Item(4, 12).Value = 0.7465491669266839
$ws.Cells.Item(4, 15).Value = 3.170708462998107
$ws.Cells.Item(5, 2).Value = 1.710224118925169
$ws.Cells.Item(5, 4).Value = 0.223169973419779
$ws.Cells.Item(5, 5).Value = 0.258834348216773
$ws.Cells.Item(5, 6).Value = 1.3202261931824
$ws.Cells.Item(5, 7).Value = 0.002453536349084012
$ws.Cells.Item(5, 10).Value = 0.3625535564687539
$ws.Cells.Item(5, 12).Value = 0.7206545757124161
$ws.Cells.Item(5, 15).Value = 3.179608155346045
$ws.Cells.Item(6, 2).Value = 1.704785686203081
$ws.Cells.Item(6, 4).Value = 0.2231660903477746
$ws.Cells.Item(6, 5).Value = 0.258582270717099
$ws.Cells.Item(6, 6).Value = 1.321190150508194
$ws.Cells.Item(6, 7).Value = 0.002453699922471575
$ws.Cells.Item(6, 10).Value = 0.3619303168644876
$ws.Cells.Item(6, 12).Value = 0.7163512270631145
$ws.Cells.Item(6, 15).Value = 3.181130337046056
$ws.Cells.Item(7, 2).Value = 1.74253681616301
$ws.Cells.Item(7, 4).Value = 0.2231993474934697
$ws.Cells.Item(7, 5).Value = 0.2603403202378018
$ws.Cells.Item(7, 6).Value = 1.314601983899919
$ws.Cells.Item(7, 7).Value = 0.00245257503013158
$ws.Cells.Item(7, 10).Value = 0.3662677012498534
$ws.Cells.Item(7, 12).Value = 0.7462001838078436
$ws.Cells.Item(7, 15).Value = 3.17082550983352
$ws.Cells.Item(8, 2).Value = 1.909173996182972
$ws.Cells.Item(8, 4).Value = 0.2235013614662975
$ws.Cells.Item(8, 5).Value = 0.2683032510595638
$ws.Cells.Item(8, 6).Value = 1.288077270736125
$ws.Cells.Item(8, 7).Value = 0.002447869969976211
$ws.Cells.Item(8, 10).Value = 0.3856870827367942
$ws.Cells.Item(8, 12).Value = 0.8773876177005207
$ws.Cells.Item(8, 15).Value = 3.13182082046626
$ws.Cells.Item(9, 2).Value = 2.236039264819567
$ws.Cells.Item(9, 4).Value = 0.2246448499111722
$ws.Cells.Item(9, 5).Value = 0.2846403734508627
$ws.Cells.Item(9, 6).Value = 1.245264686095787
$ws.Cells.Item(9, 7).Value = 0.002439567467996434
$ws.Cells.Item(9, 10).Value = 0.4247458310985763
$ws.Cells.Item(9, 12).Value = 1.132660199660563
$ws.Cells.Item(9, 15).Value = 3.078729555036318
$ws.Cells.Item(10, 2).Value = 2.476198271511691
$ws.Cells.Item(10, 4).Value = 0.2258154020381511
$ws.Cells.Item(10, 5).Value = 0.2970720624249878
$ws.Cells.Item(10, 6).Value = 1.219462927437313
$ws.Cells.Item(10, 7).Value = 0.002434025067526193
$ws.Cells.Item(10, 10).Value = 0.4540190196081255
$ws.Cells.Item(10, 12).Value = 1.318951199607511
$ws.Cells.Item(10, 15).Value = 3.054128526533418
$ws.Cells.Item(11, 2).Value = 2.585448343341056
$ws.Cells.Item(11, 4).Value = 0.2264185177725437
$ws.Cells.Item(11, 5).Value = 0.3028184770093816
$ws.Cells.Item(11, 6).Value = 1.20896161877161
$ws.Cells.Item(11, 7).Value = 0.002431623469578696
$ws.Cells.Item(11, 10).Value = 0.4674577424713959
$ws.Cells.Item(11, 12).Value = 1.403418548556601
$ws.Cells.Item(11, 15).Value = 3.046095464355176
$ws.Cells.Item(12, 2).Value = 2.626817515020036
$ws.Cells.Item(12, 4).Value = 0.2266569590859149
$ws.Cells.Item(12, 5).Value = 0.3050073988334532
$ws.Cells.Item(12, 6).Value = 1.205163548293164
$ws.Cells.Item(12, 7).Value = 0.002430731159213776
$ws.Cells.Item(12, 10).Value = 0.4725638507228496
$ws.Cells.Item(12, 12).Value = 1.435363212087168
$ws.Cells.Item(12, 15).Value = 3.043510124324939
$ws.Cells.Item(13, 2).Value = 2.617908012577914
$ws.Cells.Item(13, 4).Value = 0.2266051609076669
$ws.Cells.Item(13, 5).Value = 0.3045354061773509
$ws.Cells.Item(13, 6).Value = 1.205973576966457
$ws.Cells.Item(13, 7).Value = 0.002430922573947079
$ws.Cells.Item(13, 10).Value = 0.4714634024749387
$ws.Cells.Item(13, 12).Value = 1.428485215471994
$ws.Cells.Item(13, 15).Value = 3.044046576405634
$ws.Cells.Item(14, 2).Value = 2.588851847876128
$ws.Cells.Item(14, 4).Value = 0.2264379335515798
$ws.Cells.Item(14, 5).Value = 0.3029983041135011
$ws.Cells.Item(14, 6).Value = 1.208645566262781
$ws.Cells.Item(14, 7).Value = 0.00243154971599172
$ws.Cells.Item(14, 10).Value = 0.4678774835964816
$ws.Cells.Item(14, 12).Value = 1.406047490629135
$ws.Cells.Item(14, 15).Value = 3.04587360161068
$ws.Cells.Item(15, 2).Value = 2.571053873154995
$ws.Cells.Item(15, 4).Value = 0.2263368083705615
$ws.Cells.Item(15, 5).Value = 0.3020584551954641
$ws.Cells.Item(15, 6).Value = 1.210305512878826
$ws.Cells.Item(15, 7).Value = 0.002431936085757727
$ws.Cells.Item(15, 10).Value = 0.4656832267148445
$ws.Cells.Item(15, 12).Value = 1.39229831834615
$ws.Cells.Item(15, 15).Value = 3.047052244471274
$ws.Cells.Item(16, 2).Value = 2.469058400732308
$ws.Cells.Item(16, 4).Value = 0.2257773994763994
$ws.Cells.Item(16, 5).Value = 0.2966983354964512
$ws.Cells.Item(16, 6).Value = 1.220174153883548
$ws.Cells.Item(16, 7).Value = 0.002434184418393899
$ws.Cells.Item(16, 10).Value = 0.4531431905916179
$ws.Cells.Item(16, 12).Value = 1.313425357534072
$ws.Cells.Item(16, 15).Value = 3.054717266499694
$ws.Cells.Item(17, 2).Value = 2.406486510862521
$ws.Cells.Item(17, 4).Value = 0.2254522365830667
$ws.Cells.Item(17, 5).Value = 0.2934332616864666
$ws.Cells.Item(17, 6).Value = 1.226545399075547
$ws.Cells.Item(17, 7).Value = 0.002435594286368878
$ws.Cells.Item(17, 10).Value = 0.4454812833792232
$ws.Cells.Item(17, 12).Value = 1.264967299958698
$ws.Cells.Item(17, 15).Value = 3.060230131793702
$ws.Cells.Item(18, 2).Value = 2.37049695394046
$ws.Cells.Item(18, 4).Value = 0.2252718650878052
$ws.Cells.Item(18, 5).Value = 0.2915638731139722
$ws.Cells.Item(18, 6).Value = 1.230326264671064
$ws.Cells.Item(18, 7).Value = 0.0024364164738584
$ws.Cells.Item(18, 10).Value = 0.4410858835076539
$ws.Cells.Item(18, 12).Value = 1.237069481737137
$ws.Cells.Item(18, 15).Value = 3.063698090396883
$ws.Cells.Item(19, 2).Value = 2.358311579855183
$ws.Cells.Item(19, 4).Value = 0.2252119402733541
$ws.Cells.Item(19, 5).Value = 0.2909324137024356
$ws.Cells.Item(19, 6).Value = 1.231626349813872
$ws.Cells.Item(19, 7).Value = 0.002436696790601729
$ws.Cells.Item(19, 10).Value = 0.4395996692118729
$ws.Cells.Item(19, 12).Value = 1.227619333426219
$ws.Cells.Item(19, 15).Value = 3.064923233963611
$ws.Cells.Item(20, 2).Value = 2.413147392228211
$ws.Cells.Item(20, 4).Value = 0.2254861628444189
$ws.Cells.Item(20, 5).Value = 0.293779946458173
$ws.Cells.Item(20, 6).Value = 1.22585512731893
$ws.Cells.Item(20, 7).Value = 0.002435443037825717
$ws.Cells.Item(20, 10).Value = 0.4462957165871444
$ws.Cells.Item(20, 12).Value = 1.270128446729643
$ws.Cells.Item(20, 15).Value = 3.059612510703147
$ws.Cells.Item(21, 2).Value = 2.597386399400591
$ws.Cells.Item(21, 4).Value = 0.2264867802061303
$ws.Cells.Item(21, 5).Value = 0.303449440529505
$ws.Cells.Item(21, 6).Value = 1.207855886097576
$ws.Cells.Item(21, 7).Value = 0.00243136504525875
$ws.Cells.Item(21, 10).Value = 0.468930292566796
$ws.Cells.Item(21, 12).Value = 1.412639125857652
$ws.Cells.Item(21, 15).Value = 3.045324548105157
$ws.Cells.Item(22, 2).Value = 2.717788054440462
$ws.Cells.Item(22, 4).Value = 0.227199300565303
$ws.Cells.Item(22, 5).Value = 0.3098439806061748
$ws.Cells.Item(22, 6).Value = 1.197133424589524
$ws.Cells.Item(22, 7).Value = 0.002428799606336567
$ws.Cells.Item(22, 10).Value = 0.4838231064478009
$ws.Cells.Item(22, 12).Value = 1.505536809825117
$ws.Cells.Item(22, 15).Value = 3.038649138656695
$ws.Cells.Item(23, 2).Value = 2.653528761209429
$ws.Cells.Item(23, 4).Value = 0.2268136900386253
$ws.Cells.Item(23, 5).Value = 0.3064243100025905
$ws.Cells.Item(23, 6).Value = 1.202760678860365
$ws.Cells.Item(23, 7).Value = 0.002430159728914271
$ws.Cells.Item(23, 10).Value = 0.4758655337042796
$ws.Cells.Item(23, 12).Value = 1.455978098449521
$ws.Cells.Item(23, 15).Value = 3.041967483495938
$ws.Cells.Item(24, 2).Value = 2.410136057922159
$ws.Cells.Item(24, 4).Value = 0.2254708043120104
$ws.Cells.Item(24, 5).Value = 0.2936231861269576
$ws.Cells.Item(24, 6).Value = 1.226166831816677
$ws.Cells.Item(24, 7).Value = 0.002435511381146587
$ws.Cells.Item(24, 10).Value = 0.4459274816086349
$ws.Cells.Item(24, 12).Value = 1.267795212534622
$ws.Cells.Item(24, 15).Value = 3.059890807428019
$ws.Cells.Item(25, 2).Value = 2.147608612048373
$ws.Cells.Item(25, 4).Value = 0.2242770781105321
$ws.Cells.Item(25, 5).Value = 0.2801446540664614
$ws.Cells.Item(25, 6).Value = 1.255857038352787
$ws.Cells.Item(25, 7).Value = 0.002441715195963146
$ws.Cells.Item(25, 10).Value = 0.4140767872051896
$ws.Cells.Item(25, 12).Value = 1.063819820018466
$ws.Cells.Item(25, 15).Value = 3.090572688611672
